$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.750.22"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.868.92"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'0.7213"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").Value = "'240.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.3127"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.07129"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'24.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.08123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").Value = "1.909.65"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "'0.7416"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "'5.343"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'92.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "29.790.55"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'5.989"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "'246.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "'13.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'0.000007789"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "2.141.42"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "'7.718"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "'0.1521"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("D26").Value = "'9.216"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'163.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "'18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'2.006"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "'1.443"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "'4.513"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("D32").Value = "'1.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'4.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").Value = "'0.05363"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "'1.225"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").Value = "'0.7359"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "'2.694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'0.01924"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "'2.737"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "'0.4462"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "'0.8832"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "'5.966"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "'71.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.039.43"
$ws.Range("E45").Value = "  -6.55%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'103.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'7.457"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").Value = "'1.815"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "'9.549"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "2.032.08"
$ws.Range("E51").Value = "  +0.52%  "
